# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# text values on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 17:40:15"
$wsZhCn.Range("E3").Value = "2016-03-20 17:40:15"
$wsZhCn.Range("H2").Value = "2016-03-20 17:40:56"
$wsZhCn.Range("H3").Value = "2016-03-20 17:40:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 17:40:24"
$wsDeDe.Range("E3").Value = "2016-03-20 17:40:24"
$wsDeDe.Range("H2").Value = "2016-03-20 17:41:10"
$wsDeDe.Range("H3").Value = "2016-03-20 17:41:10"
